# Update "countries & provincias Spain" data sheet.
#
# The underlying change re-orders a couple of shared-string rows
# ("Tenerife" moves up next to "Gran Canaria"/"Soria"; "La Gomera" moves
# up next to "Menorca"/"Arroyo de la Luz") and refreshes the timestamp
# footer plus a batch of the day's case counts. Because inserting a row
# shifts every following row's data up by one, the net effect (as seen
# by anyone reading the sheet) is that a contiguous block of rows ends
# up showing new city labels together with the numbers that used to
# belong to the row above, followed by a handful of rows whose counts
# were simply refreshed with newer figures.
#
# We reproduce that end state directly, cell by cell, through the
# Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Mayo de 2020 a las 16:35"

# "Tenerife" now sits right after "Gran Canaria" (row 34), pushing the
# following provinces down one row each (through "Pontevedra", row 42).
$ws.Range("A34").Value = "Tenerife"
$ws.Range("B34").Value = 2280
$ws.Range("C34").Value = 1506
$ws.Range("D34").Value = 623
$ws.Range("E34").Value = 151

$ws.Range("A35").Value = "Soria"
$ws.Range("B35").Value = 2278
$ws.Range("C35").Value = 393
$ws.Range("D35").Value = 1766
$ws.Range("E35").Value = 119

$ws.Range("A36").Value = "Cantabria"
$ws.Range("B36").Value = 2246
$ws.Range("C36").Value = 1981
$ws.Range("D36").Value = 62
$ws.Range("E36").Value = 203

$ws.Range("A37").Value = "Caceres"
$ws.Range("B37").Value = 1973
$ws.Range("C37").Value = 1505
$ws.Range("D37").Value = 66
$ws.Range("E37").Value = 402

$ws.Range("A38").Value = "A Coruña"
$ws.Range("B38").Value = 1969
$ws.Range("C38").Value = 333
$ws.Range("D38").Value = 1788
$ws.Range("E38").Value = 67

$ws.Range("A39").Value = "Avila"
$ws.Range("B39").Value = 1917
$ws.Range("C39").Value = 618
$ws.Range("D39").Value = 1166
$ws.Range("E39").Value = 133

$ws.Range("A40").Value = "Jaen"
$ws.Range("B40").Value = 1751
$ws.Range("C40").Value = 1121
$ws.Range("D40").Value = 457
$ws.Range("E40").Value = 173

$ws.Range("A41").Value = "Cordoba"
$ws.Range("B41").Value = 1682
$ws.Range("C41").Value = 1331
$ws.Range("D41").Value = 246
$ws.Range("E41").Value = 105

$ws.Range("A42").Value = "Pontevedra"
$ws.Range("B42").Value = 1536
$ws.Range("C42").Value = 333
$ws.Range("D42").Value = 1411
$ws.Range("E42").Value = 30

# (row 43 "Murcia" onward through row 59 already had and keeps their
# current values - untouched)

# Refreshed case counts for "La Palma" / "Lanzarote" (rows unaffected by
# the row re-ordering, just newer figures).
$ws.Range("B60").Value = 95
$ws.Range("C60").Value = 68
$ws.Range("D60").Value = 22
$ws.Range("E60").Value = 5

$ws.Range("B61").Value = 84
$ws.Range("C61").Value = 71
$ws.Range("D61").Value = 7
$ws.Range("E61").Value = 6

# Refreshed figures for "Fuerteventura"
$ws.Range("C63").Value = 42
$ws.Range("D63").Value = 0

# "La Gomera" now sits right before "Arroyo de la Luz".
$ws.Range("A66").Value = "La Gomera"
$ws.Range("B66").Value = 8
$ws.Range("C66").Value = 8
$ws.Range("D66").Value = 0

$ws.Range("A67").Value = "Arroyo de la Luz"
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 7

# Refreshed figures for "El Hierro"
$ws.Range("B68").Value = 3
$ws.Range("C68").Value = 3
